$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = -0.1140116980535687
$ws.Range("C2").Value = 1.134027854149689
$ws.Range("D2").Value = 0.5293419069744123
$ws.Range("F2").Value = 1.285355063281015
$ws.Range("G2").Value = 1.444240868935941
$ws.Range("H2").Value = 0.3125477243455349
$ws.Range("I2").Value = -0.5548516401491572
$ws.Range("B3").Value = 0.52162968558209
$ws.Range("C3").Value = 0.6031613697014773
$ws.Range("E3").Value = 1.4345885005637
$ws.Range("F3").Value = 1.547161051919687
$ws.Range("G3").Value = 0.4466123193550214
$ws.Range("H3").Value = -0.4809009009009035
$ws.Range("B4").Value = -0.5273222390317647
$ws.Range("D4").Value = -0.005411499436300538
$ws.Range("E4").Value = 0.7074136955291455
$ws.Range("F4").Value = 0.3090111142775654
$ws.Range("G4").Value = -0.4209009009009035
$ws.Range("H4").Value = 1.271278195488724
$ws.Range("I4").Value = 0.6220589070934801
$ws.Range("J4").Value = 0.8193626707131937
$ws.Range("C5").Value = 0.004588500563699471
$ws.Range("D5").Value = 0.7769453239451327
$ws.Range("E5").Value = 0.3807516038255184
$ws.Range("F5").Value = -0.3988174318051805
$ws.Range("G5").Value = 1.180575399060113
$ws.Range("H5").Value = 0.5077756378413858
$ws.Range("I5").Value = 0.6787354624760854
$ws.Range("B6").Value = 0.1867105431484077
$ws.Range("C6").Value = 0.7722779613039183
$ws.Range("D6").Value = 0.2843558036209001
$ws.Range("E6").Value = -0.3395270793998535
$ws.Range("F6").Value = 1.287674180139613
$ws.Range("G6").Value = 0.5677696640183092
$ws.Range("H6").Value = 0.6675525209407469
$ws.Range("B7").Value = 0.8280615287077258
$ws.Range("C7").Value = 0.009011114277565485
$ws.Range("D7").Value = -0.4809009009009035
$ws.Range("E7").Value = 1.261278195488724
$ws.Range("F7").Value = 0.6020589070934801
$ws.Range("G7").Value = 0.7093626707131937
$ws.Range("B8").Value = -0.3509888857224345
$ws.Range("C8").Value = -0.5809009009009034
$ws.Range("D8").Value = 1.241278195488724
$ws.Range("E8").Value = 0.7032066015985186
$ws.Range("F8").Value = 0.8825435650546625
$ws.Range("H8").Value = -0.1837499599687909
$ws.Range("I8").Value = 0.3538112635634896
$ws.Range("J8").Value = 0.1485032540344368
$ws.Range("B9").Value = -0.8569009009009034
$ws.Range("C9").Value = 1.055249540297363
$ws.Range("D9").Value = 0.5422953553339147
$ws.Range("E9").Value = 0.8185372705227552
$ws.Range("G9").Value = -0.115966644056028
$ws.Range("H9").Value = 0.3949757811758082
$ws.Range("I9").Value = 0.1734537503564907
$ws.Range("B10").Value = 0.6787681550252201
$ws.Range("C10").Value = 0.3564167683417685
$ws.Range("D10").Value = 0.8907192198060601
$ws.Range("F10").Value = -0.0717095860245372
$ws.Range("G10").Value = 0.4455085207778247
$ws.Range("H10").Value = 0.2422520263583712
$ws.Range("B11").Value = 0.4003114290510102
$ws.Range("C11").Value = 0.8093626707131937
$ws.Range("E11").Value = -0.05551312649165269
$ws.Range("F11").Value = 0.4849757811758083
$ws.Range("G11").Value = 0.2534537503564908
$ws.Range("B12").Value = 0.33218714548825
$ws.Range("D12").Value = -0.2025131264916527
$ws.Range("E12").Value = 0.5149757811758082
$ws.Range("F12").Value = 0.381103329907261
$ws.Range("C13").Value = -0.3157490177063643
$ws.Range("D13").Value = 0.2454495103302783
$ws.Range("E13").Value = 0.04235042473292953
$ws.Range("B14").Value = -0.2760759690770392
$ws.Range("C14").Value = 0.2958161606567877
$ws.Range("D14").Value = 0.07961008106920435
$ws.Range("B15").Value = 0.1421546153588515
$ws.Range("C15").Value = 0.02893023050567838
$ws.Range("B16").Value = 0.02940328597706714
